$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value (mirrors the authoritative diff)
$updates = @{
    "D2" = '26.772.84'
    "E2" = '  +0.04%  '
    "D3" = '1.642.32'
    "E3" = '  -0.52%  '
    "E4" = '  +0.57%  '
    "D5" = '216.83'
    "E5" = '  +0.34%  '
    "E6" = '  -0.81%  '
    "E7" = '  +0.52%  '
    "E8" = '  -1.06%  '
    "E9" = '  -0.78%  '
    "D10" = '19.20'
    "E10" = '  -1.34%  '
    "D11" = '0.0840'
    "E11" = '  -0.85%  '
    "D12" = '1.869.13'
    "E12" = '  -0.57%  '
    "D13" = '1.641.40'
    "E13" = '  -0.81%  '
    "E14" = '  -1.55%  '
    "E15" = '  -1.88%  '
    "D16" = '64.59'
    "E16" = '  -3.45%  '
    "D17" = '26.783.11'
    "E17" = '  -0.09%  '
    "E18" = '  -2.62%  '
    "D19" = '214.32'
    "E19" = '  -3.18%  '
    "E20" = '  +0.57%  '
    "E21" = '  -1.64%  '
    "E22" = '  +11.80%  '
    "D23" = '6.28'
    "E23" = '  -1.27%  '
    "D24" = '9.36'
    "E24" = '  -2.59%  '
    "D25" = '144.91'
    "E25" = '  -1.41%  '
    "E26" = '  +0.88%  '
    "D27" = '0.119'
    "E27" = '  -2.39%  '
    "E28" = '  -0.45%  '
    "D29" = '15.68'
    "E29" = '  -1.86%  '
    "E30" = '  -1.62%  '
    "E31" = '  -0.02%  '
    "E32" = '  -3.51%  '
    "D33" = '3.00'
    "E33" = '  -2.36%  '
    "D34" = '1.291.27'
    "E34" = '  -0.32%  '
    "E35" = '  -2.21%  '
    "E36" = '  +1.07%  '
    "E37" = '  -4.16%  '
    "D38" = '0.540'
    "E38" = '  +1.90%  '
    "D39" = '0.825'
    "E39" = '  -1.30%  '
    "E40" = '  +0.50%  '
    "D41" = '0.810'
    "E41" = '  -0.95%  '
    "D42" = '2.24'
    "E42" = '  -0.25%  '
    "D44" = '1.795.09'
    "E45" = '  -2.51%  '
    "D46" = '59.94'
    "E46" = '  -0.09%  '
    "E47" = '  -1.17%  '
    "E48" = '  -2.00%  '
    "D49" = '0.0521'
    "E49" = '  +0.57%  '
    "D50" = '7.68'
    "E50" = '  -2.15%  '
    "D51" = '0.0977'
    "E51" = '  -0.60%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage so numeric-looking strings (e.g. "19.20", "0.0840")
    # keep their exact digits/trailing zeros instead of being parsed as numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    # Restore the default (unstyled) cell style so no stray formatting is introduced -
    # these data cells carry no style index in the source workbook.
    $cell.Style = "Normal"
}
